# Add season record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: AD1 = "Wins", AE1 = "Losses", AF1 = "Ties"
# Copy the formatting from the existing header cell (A1) onto the new
# header cells so they share the same bold/bordered/centered style.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2 through 48: AD = 78 (Wins), AE = 84 (Losses), AF = 0 (Ties)
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 78
    $ws.Cells.Item($r, 31).Value = 84
    $ws.Cells.Item($r, 32).Value = 0
}
